$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# -----------------------------------------------------------------
# 1. Update the "web" column (AE) on the '#system' sheet to insert
#    the two new commands in their alphabetically-sorted position:
#       assertSelectOptionsAbsent(locator,options)
#       assertSelectOptionsPresent(locator,options)
#    This pushes every entry from AE39 downward by two rows, so the
#    list now spans AE2:AE158 (was AE2:AE156).
# -----------------------------------------------------------------
$webValues = @(
    "assertSelectOptionsAbsent(locator,options)",
    "assertSelectOptionsPresent(locator,options)",
    "assertSingleSelect(locator)",
    "assertTable(locator,row,column,text)",
    "assertText(locator,text)",
    "assertTextContains(locator,text)",
    "assertTextCount(locator,text,count)",
    "assertTextList(locator,list,ignoreOrder)",
    "assertTextMatches(text,minMatch,scrollTo)",
    "assertTextNotContain(locator,text)",
    "assertTextNotPresent(text)",
    "assertTextOrder(locator,descending)",
    "assertTextPresent(text)",
    "assertTitle(text)",
    "assertValue(locator,value)",
    "assertValueOrder(locator,descending)",
    "assertVisible(locator)",
    "checkAll(locator,waitMs)",
    "clearLocalStorage()",
    "click(locator)",
    "clickAll(locator)",
    "clickAndWait(locator,waitMs)",
    "clickByLabel(label)",
    "clickByLabelAndWait(label,waitMs)",
    "clickOffset(locator,x,y)",
    "clickWithKeys(locator,keys)",
    "close()",
    "closeAll()",
    "deselect(locator,text)",
    "deselectMulti(locator,array)",
    "dismissInvalidCert()",
    "dismissInvalidCertPopup()",
    "doubleClick(locator)",
    "doubleClickAndWait(locator,waitMs)",
    "doubleClickByLabel(label)",
    "doubleClickByLabelAndWait(label,waitMs)",
    "dragAndDrop(fromLocator,toLocator)",
    "dragTo(fromLocator,xOffset,yOffset)",
    "editLocalStorage(key,value)",
    "executeScript(var,script)",
    "focus(locator)",
    "goBack()",
    "goBackAndWait()",
    "maximizeWindow()",
    "mouseOver(locator)",
    "moveTo(x,y)",
    "open(url)",
    "openAndWait(url,waitMs)",
    "openHttpBasic(url,username,password)",
    "openIgnoreTimeout(url)",
    "openInTab(name,url)",
    "refresh()",
    "refreshAndWait()",
    "resizeWindow(width,height)",
    "rightClick(locator)",
    "saveAllWindowIds(var)",
    "saveAllWindowNames(var)",
    "saveAttribute(var,locator,attrName)",
    "saveAttributeList(var,locator,attrName)",
    "saveBrowserVersion(var)",
    "saveCount(var,locator)",
    "saveCssValue(var,locator,property)",
    "saveDivsAsCsv(headers,rows,cells,nextPage,file)",
    "saveElement(var,locator)",
    "saveElements(var,locator)",
    "saveInfiniteDivsAsCsv(config,file)",
    "saveInfiniteTableAsCsv(config,file)",
    "saveLocalStorage(var,key)",
    "saveLocation(var)",
    "savePageAs(var,sessionIdName,url)",
    "savePageAsFile(sessionIdName,url,file)",
    "saveSelectedText(var,locator)",
    "saveSelectedValue(var,locator)",
    "saveTableAsCsv(locator,nextPageLocator,file)",
    "saveText(var,locator)",
    "saveTextArray(var,locator)",
    "saveTextSubstringAfter(var,locator,delim)",
    "saveTextSubstringBefore(var,locator,delim)",
    "saveTextSubstringBetween(var,locator,start,end)",
    "saveTitle(var)",
    "saveValue(var,locator)",
    "saveValues(var,locator)",
    "screenshot(file,locator,removeFixed)",
    "screenshotInFull(file,timeout,removeFixed)",
    "scrollPage(xOffset,yOffset)",
    "scrollTo(locator)",
    "select(locator,text)",
    "selectAllOptions(locator)",
    "selectDropdown(locator,optLocator,optText)",
    "selectFrame(locator)",
    "selectMulti(locator,array)",
    "selectMultiByValue(locator,array)",
    "selectMultiOptions(locator)",
    "selectText(locator)",
    "selectWindow(winId)",
    "selectWindowAndWait(winId,waitMs)",
    "selectWindowByIndex(index)",
    "selectWindowByIndexAndWait(index,waitMs)",
    "switchBrowser(profile,config)",
    "toggleSelections(locator)",
    "type(locator,value)",
    "typeKeys(locator,value)",
    "uncheckAll(locator,waitMs)",
    "unselectAllText()",
    "updateAttribute(locator,attrName,value)",
    "upload(fieldLocator,file)",
    "verifyContainText(locator,text)",
    "verifyText(locator,text)",
    "wait(waitMs)",
    "waitForElementPresent(locator,waitMs)",
    "waitForElementTextPresent(locator,text)",
    "waitForElementsPresent(locators)",
    "waitForPopUp(winId,waitMs)",
    "waitForTextPresent(text)",
    "waitForTitle(text)",
    "waitUntilDisabled(locator,waitMs)",
    "waitUntilEnabled(locator,waitMs)",
    "waitUntilHidden(locator,waitMs)",
    "waitUntilVisible(locator,waitMs)",
    "waitWhileElementNotPresent(locator,waitMs)"
)

for ($i = 0; $i -lt $webValues.Length; $i++) {
    $row = 39 + $i
    $ws.Cells.Item($row, 31).Value = $webValues[$i]
}

# -----------------------------------------------------------------
# 2. Update the "web" defined name so it refers to the new, larger
#    range (was $AE$2:$AE$156, now $AE$2:$AE$158).
# -----------------------------------------------------------------
$wb.Names.Item("web").RefersTo = "='#system'!`$AE`$2:`$AE`$158"

# -----------------------------------------------------------------
# 3. Touch column AM formatting so the sheet's recorded dimension
#    keeps its original right-hand edge (column AM) after the two
#    new rows are appended, matching "A1:AM158".
# -----------------------------------------------------------------
$ws.Range("AM1:AM158").Font.Bold = $false

Write-Output "done"
